$d = $word.ActiveDocument
$rng = $d.Content
$rng.Collapse(1)
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid wp14"><w:body><w:p w14:paraId="483F45E7" w14:textId="14DB1693" w:rsidR="002050DA" w:rsidRDefault="002050DA" w:rsidP="002050DA"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>DESIGN RATIONALE: FIT2099 SSB ASSIGNMENT 1</w:t></w:r></w:p><w:p w14:paraId="3E698ADB" w14:textId="77777777" w:rsidR="002050DA" w:rsidRDefault="002050DA" w:rsidP="002050DA"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">GROUP NAME </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>JAVAPLUSPLUS</w:t></w:r></w:p><w:p w14:paraId="7200C886" w14:textId="77777777" w:rsidR="002050DA" w:rsidRDefault="002050DA" w:rsidP="002050DA"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">GROUP MEMBERS </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t xml:space="preserve">MOHAMED SHAKEEL MOHAMED RAFI </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>28021452</w:t></w:r></w:p><w:p w14:paraId="2319BC23" w14:textId="77777777" w:rsidR="002050DA" w:rsidRDefault="002050DA" w:rsidP="002050DA"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t xml:space="preserve">      </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>MATTI HADDAD</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>29708966</w:t></w:r></w:p><w:p w14:paraId="2D533D25" w14:textId="77777777" w:rsidR="002050DA" w:rsidRDefault="002050DA" w:rsidP="002050DA"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>KERRY YUE SONG ZHENG</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>28794346</w:t></w:r></w:p><w:p w14:paraId="440D2B34" w14:textId="77777777" w:rsidR="002050DA" w:rsidRDefault="002050DA" w:rsidP="007E32B7"/><w:p w14:paraId="2A920C54" w14:textId="77777777" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"><w:r><w:t>This document relays the design rationale and thought process for the tasks that were assigned.</w:t></w:r></w:p><w:p w14:paraId="3FEFA937" w14:textId="77777777" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"><w:r><w:t>This includes:</w:t></w:r></w:p><w:p w14:paraId="4AE0309A" w14:textId="36C567A1" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Leave Affordance</w:t></w:r></w:p><w:p w14:paraId="154C9A3A" w14:textId="5F89DF5F" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Give Affordance</w:t></w:r></w:p><w:p w14:paraId="728AE4B0" w14:textId="28EA99B9" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Wand and Spell </w:t></w:r><w:r w:rsidR="00942754"><w:t>Class</w:t></w:r><w:r><w:t xml:space="preserve"> Implementation</w:t></w:r></w:p><w:p w14:paraId="3534D72C" w14:textId="167231DC" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="002050DA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Dementor</w:t></w:r><w:r w:rsidR="002050DA"><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00942754"><w:t>SubClass</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00942754"><w:t xml:space="preserve"> of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00942754"><w:t>HP</w:t></w:r><w:r w:rsidR="002050DA"><w:t>Actor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="002050DA"><w:t>) Implementation</w:t></w:r></w:p><w:p w14:paraId="04E43BB8" w14:textId="77777777" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"/><w:p w14:paraId="42D3AC8E" w14:textId="5B2922B3" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="002050DA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Leave Affordance</w:t></w:r></w:p><w:p w14:paraId="16EB0F0D" w14:textId="77777777" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"><w:r><w:t xml:space="preserve">After an item is ''taken'', it needs to have the option to be ''left'' as well. </w:t></w:r></w:p><w:p w14:paraId="52AD2E7A" w14:textId="77777777" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"><w:r><w:t>Hence, leave needs to be a type of affordance, which in the code, will eventually give the actor the action of leaving it.</w:t></w:r></w:p><w:p w14:paraId="2D8275DE" w14:textId="6239AE6D" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"><w:r><w:t>Once an item has been left by an actor, the Take affordance needs to be reassigned to the item so that an actor could take it again.</w:t></w:r></w:p><w:p w14:paraId="5140420C" w14:textId="2AC5107F" w:rsidR="002050DA" w:rsidRDefault="00942754" w:rsidP="002050DA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Give</w:t></w:r><w:r w:rsidR="002050DA"><w:t xml:space="preserve"> Affordance</w:t></w:r></w:p><w:p w14:paraId="651690C8" w14:textId="28371F7F" w:rsidR="002A336A" w:rsidRDefault="002A336A" w:rsidP="002A336A"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/></w:pPr><w:r><w:t>When an item picked</w:t></w:r><w:r><w:t xml:space="preserve"> by an actor</w:t></w:r><w:r><w:t xml:space="preserve">, the actor has the option of </w:t></w:r><w:r><w:t>‘</w:t></w:r><w:r><w:t>giving</w:t></w:r><w:r><w:t>’</w:t></w:r><w:r><w:t xml:space="preserve"> the item to another actor</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">if and only if </w:t></w:r><w:r><w:t>both</w:t></w:r><w:r><w:t xml:space="preserve"> actors are in</w:t></w:r><w:r><w:t xml:space="preserve"> the same location</w:t></w:r><w:r><w:t xml:space="preserve"> and both are from one team</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p w14:paraId="65C0153D" w14:textId="6ECAFFF9" w:rsidR="002A336A" w:rsidRDefault="002A336A" w:rsidP="002A336A"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/></w:pPr><w:r><w:t>Thus</w:t></w:r><w:r><w:t>, a give class</w:t></w:r><w:r><w:t xml:space="preserve"> which deals with entities,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">is needed </w:t></w:r><w:r><w:t xml:space="preserve">to perform such </w:t></w:r><w:r><w:t>action</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/></w:pPr><w:r><w:t>The Give class</w:t></w:r><w:r><w:t xml:space="preserve"> is a subclass of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>HPAffordance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> that </w:t></w:r><w:r><w:t>depends on:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>MessageRender</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to display a message to the user, a form of player feedback.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>HPEntetityInterface</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as</w:t></w:r><w:r><w:t xml:space="preserve"> it deals with entities and therefore</w:t></w:r><w:r><w:t xml:space="preserve"> it needs to be managed by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>EntityManager</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>HPActor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as actors are the ones who initiate the action</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>HPAction</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as </w:t></w:r><w:r><w:t>to check if the action can be performed.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">It also </w:t></w:r><w:r><w:t>implements</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>HPActionInterface</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>as it</w:t></w:r><w:r><w:t xml:space="preserve"> will need to</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>implement some method from it.</w:t></w:r></w:p><w:p w14:paraId="47F7F6A2" w14:textId="5642FD6C" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="00D6006F"/><w:p w14:paraId="63010C12" w14:textId="77777777" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="00D6006F"/><w:p w14:paraId="76DC7DD3" w14:textId="1E371C24" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="002050DA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Wand and Spell Implementation</w:t></w:r></w:p><w:p w14:paraId="7774F75B" w14:textId="77777777" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"><w:r><w:t>For a Spell to be Casted, there must be a Wand class and a Cast Action / Affordance.</w:t></w:r></w:p><w:p w14:paraId="5D425259" w14:textId="77777777" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="007E32B7"><w:r><w:t>Furthermore, Cast may target both entities and actors, which will require further distinction.</w:t></w:r></w:p><w:p w14:paraId="0B6A74BB" w14:textId="3A3907B3" w:rsidR="007E32B7" w:rsidRDefault="002050DA" w:rsidP="002050DA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Wand and Cast Class</w:t></w:r></w:p><w:p w14:paraId="51A6878F" w14:textId="77777777" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="00151560"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>A Wand, like a dagger or any other item in the game, can be picked up by any Actor. The existence of a Wand allows an Actor to Cast Spells onto other Actors, or itself.</w:t></w:r></w:p><w:p w14:paraId="38E92CDA" w14:textId="77777777" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="00D6006F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/></w:pPr></w:p><w:p w14:paraId="77106306" w14:textId="646A8454" w:rsidR="002050DA" w:rsidRDefault="002050DA" w:rsidP="002050DA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Spell Abstraction</w:t></w:r></w:p><w:p w14:paraId="3B8E0EEB" w14:textId="1C6BF8E5" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="00151560"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">Every spell needs to have a common target attribute, and a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>common  method</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> which denotes the spell''s function.</w:t></w:r></w:p><w:p w14:paraId="5C801EF6" w14:textId="32534AD5" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="001063B2"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>Hence, abstraction of the Spell class allows subclasses to share that, and it makes it easier to Cast any subclass of Spell.</w:t></w:r></w:p><w:p w14:paraId="3A3E6FE5" w14:textId="0D4C8338" w:rsidR="002050DA" w:rsidRDefault="002050DA" w:rsidP="002050DA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Cast and Actor’s Known Spells</w:t></w:r></w:p><w:p w14:paraId="37FE3EFB" w14:textId="77777777" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="00151560"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">Casting a spell requires a new action - Cast. This is different from Attack as it requires a Spell to be an input as well. </w:t></w:r></w:p><w:p w14:paraId="2FC2FD06" w14:textId="49AA0FAB" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="00151560"><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:r><w:t xml:space="preserve">It will then execute the Spell''s effect if the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Actorknows</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the Spell.</w:t></w:r></w:p><w:p w14:paraId="6E9ADCB7" w14:textId="39F31AA8" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="00151560"><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:r><w:t>Every actor will have a set of its known Spells</w:t></w:r></w:p><w:p w14:paraId="3CC78E72" w14:textId="77777777" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="00D6006F"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p w14:paraId="4DCC07E5" w14:textId="3FD56AFD" w:rsidR="002050DA" w:rsidRDefault="002050DA" w:rsidP="002050DA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Targeting Items</w:t></w:r></w:p><w:p w14:paraId="128E1972" w14:textId="77777777" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="00151560"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">In order to allow Spells to be Casted on items, the interface </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>HPEntityInterface</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is used in Cast and Spell to target both subclasses.</w:t></w:r></w:p><w:p w14:paraId="5A7A0C7B" w14:textId="77777777" w:rsidR="00D6006F" w:rsidRDefault="00D6006F" w:rsidP="00D6006F"><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p w14:paraId="3EBC08BC" w14:textId="4DB1EB4A" w:rsidR="007E32B7" w:rsidRDefault="007E32B7" w:rsidP="006B576E"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">Additional checks </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>have to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> be implemented for individual spells to ensure that the target is the intended class type.</w:t></w:r></w:p><w:p w14:paraId="3D9C5425" w14:textId="21A4A40E" w:rsidR="006B576E" w:rsidRDefault="00AF014A" w:rsidP="00AF014A"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3525"/></w:tabs><w:ind w:left="360"/></w:pPr><w:r><w:tab/></w:r></w:p><w:p w14:paraId="443669E4" w14:textId="77777777" w:rsidR="00AF014A" w:rsidRDefault="00AF014A" w:rsidP="00AF014A"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3525"/></w:tabs><w:ind w:left="360"/></w:pPr></w:p><w:p w14:paraId="1EAC21AE" w14:textId="77777777" w:rsidR="006B576E" w:rsidRDefault="006B576E" w:rsidP="006B576E"><w:pPr><w:ind w:left="360"/></w:pPr></w:p><w:sectPr w:rsidR="006B576E"><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/numbering.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.numbering+xml"><pkg:xmlData><?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:numbering xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid wp14"><w:abstractNum w:abstractNumId="0" w15:restartNumberingAfterBreak="0"><w:nsid w:val="33BB4CCA"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="6A8AC46C"/><w:lvl w:ilvl="0" w:tplc="0C090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1080" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="0C090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1800" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="0C090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2520" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="0C090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3240" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="0C090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3960" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="0C090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4680" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="0C090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5400" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="0C090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6120" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="0C090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6840" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="1" w15:restartNumberingAfterBreak="0"><w:nsid w:val="4E240CD5"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="127EB50C"/><w:lvl w:ilvl="0" w:tplc="3ECEAF50"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:suff w:val="space"/><w:lvlText w:val="%1."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="720"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="44090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%2."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1080" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="2" w:tplc="4409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%3."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="1800" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="3" w:tplc="4409000F" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%4."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2520" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="4" w:tplc="44090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%5."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3240" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="5" w:tplc="4409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%6."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="3960" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="6" w:tplc="4409000F" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%7."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4680" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="7" w:tplc="44090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%8."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5400" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="8" w:tplc="4409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%9."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="6120" w:hanging="180"/></w:pPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="2" w15:restartNumberingAfterBreak="0"><w:nsid w:val="4FAD6340"/><w:multiLevelType w:val="multilevel"/><w:tmpl w:val="75C6B55C"/><w:lvl w:ilvl="0"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1080" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1800" w:hanging="720"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="720"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4.%5."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="1080"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4.%5.%6."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3240" w:hanging="1080"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4.%5.%6.%7."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3960" w:hanging="1440"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4.%5.%6.%7.%8."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="1440"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4.%5.%6.%7.%8.%9."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="1800"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="3" w15:restartNumberingAfterBreak="0"><w:nsid w:val="4FDD0DB6"/><w:multiLevelType w:val="multilevel"/><w:tmpl w:val="6E24CD48"/><w:lvl w:ilvl="0"><w:start w:val="3"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="360" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1.%2"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1068" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1.%2.%3"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="720"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1.%2.%3.%4"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="720"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1.%2.%3.%4.%5"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1080" w:hanging="1080"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1.%2.%3.%4.%5.%6"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1080" w:hanging="1080"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1.%2.%3.%4.%5.%6.%7"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="1440"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1.%2.%3.%4.%5.%6.%7.%8"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="1440"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1.%2.%3.%4.%5.%6.%7.%8.%9"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="1440"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="4" w15:restartNumberingAfterBreak="0"><w:nsid w:val="54242BB8"/><w:multiLevelType w:val="multilevel"/><w:tmpl w:val="E3EEDE72"/><w:lvl w:ilvl="0"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="360" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="720"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1800" w:hanging="720"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4.%5"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2520" w:hanging="1080"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4.%5.%6"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="1080"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4.%5.%6.%7"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="1440"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4.%5.%6.%7.%8"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3960" w:hanging="1440"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:isLgl/><w:lvlText w:val="%1.%2.%3.%4.%5.%6.%7.%8.%9"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="1440"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="5" w15:restartNumberingAfterBreak="0"><w:nsid w:val="66AA4877"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="ECB8F194"/><w:lvl w:ilvl="0" w:tplc="94C02990"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="44090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%2."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="2" w:tplc="4409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%3."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="2160" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="3" w:tplc="4409000F" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%4."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="4" w:tplc="44090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%5."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="5" w:tplc="4409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%6."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="4320" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="6" w:tplc="4409000F" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%7."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="7" w:tplc="44090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%8."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="8" w:tplc="4409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%9."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="6480" w:hanging="180"/></w:pPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="6" w15:restartNumberingAfterBreak="0"><w:nsid w:val="75225C45"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="C3C4EC1A"/><w:lvl w:ilvl="0" w:tplc="44090003"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1800" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="A25C2244"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2520" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="44090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3240" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="44090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3960" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="44090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4680" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="44090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5400" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="44090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6120" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="44090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6840" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="44090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="7560" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:num w:numId="1"><w:abstractNumId w:val="2"/></w:num><w:num w:numId="2"><w:abstractNumId w:val="3"/></w:num><w:num w:numId="3"><w:abstractNumId w:val="6"/></w:num><w:num w:numId="4"><w:abstractNumId w:val="1"/></w:num><w:num w:numId="5"><w:abstractNumId w:val="4"/></w:num><w:num w:numId="6"><w:abstractNumId w:val="5"/></w:num><w:num w:numId="7"><w:abstractNumId w:val="0"/></w:num></w:numbering></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)
Write-Host "done"
